# Regenerate save_data to use K instead of Strike#, recompute std/mean,
# and write the recalculated K values (s_vals) into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 7
